$d = $word.ActiveDocument

# The last paragraph in the body ends with "悬都是缘". We need to append a
# new run "。星光和荆棘都敢追。" (sz/szCs 18, eastAsia hint) right after it,
# still inside that same paragraph (before the end-of-paragraph mark).

$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)

# Range that covers the paragraph's text only (excludes the pilcrow) so an
# InsertAfter appends inside the paragraph instead of creating a new one.
$r = $lastPara.Range
$r.End = $r.End - 1
$r.Collapse(0)   # wdCollapseEnd

$r.InsertAfter("。星光和荆棘都敢追。")
$r.Font.Size = 9
